# Auto-generated edit script applying the commit diff
$wb = $excel.ActiveWorkbook

# ---- Section_A timetable updates ----
$wsA = $wb.Worksheets.Item("Section_A")
$aTop = New-Object 'object[,]' 2,5
$aTop[0,0] = "Free"
$aTop[0,1] = "Free"
$aTop[0,2] = "Free"
$aTop[0,3] = "Free"
$aTop[0,4] = "CS303"
$aTop[1,0] = "CS309"
$aTop[1,1] = "Free"
$aTop[1,2] = "CS309"
$aTop[1,3] = "Free"
$aTop[1,4] = "Free"
$wsA.Range("B2:F3").Value = $aTop
$aBot = New-Object 'object[,]' 4,5
$aBot[0,0] = "CS303"
$aBot[0,1] = "CS304"
$aBot[0,2] = "CS304"
$aBot[0,3] = "Free"
$aBot[0,4] = "Free"
$aBot[1,0] = "Free"
$aBot[1,1] = "Free"
$aBot[1,2] = "Free"
$aBot[1,3] = "CS303 (Tutorial)"
$aBot[1,4] = "Free"
$aBot[2,0] = "CS304"
$aBot[2,1] = "Free"
$aBot[2,2] = "CS303"
$aBot[2,3] = "Free"
$aBot[2,4] = "CS309"
$aBot[3,0] = "CS309 (Tutorial)"
$aBot[3,1] = "Free"
$aBot[3,2] = "Free"
$aBot[3,3] = "Free"
$aBot[3,4] = "CS304 (Tutorial)"
$wsA.Range("B5:F8").Value = $aBot

# ---- Section_B timetable updates ----
$wsB = $wb.Worksheets.Item("Section_B")
$bTop = New-Object 'object[,]' 2,5
$bTop[0,0] = "Free"
$bTop[0,1] = "CS303"
$bTop[0,2] = "Free"
$bTop[0,3] = "Free"
$bTop[0,4] = "CS304"
$bTop[1,0] = "CS304"
$bTop[1,1] = "Free"
$bTop[1,2] = "CS303"
$bTop[1,3] = "Free"
$bTop[1,4] = "CS309"
$wsB.Range("B2:F3").Value = $bTop
$bBot = New-Object 'object[,]' 4,5
$bBot[0,0] = "Free"
$bBot[0,1] = "Free"
$bBot[0,2] = "Free"
$bBot[0,3] = "CS303"
$bBot[0,4] = "Free"
$bBot[1,0] = "CS309 (Tutorial)"
$bBot[1,1] = "CS303 (Tutorial)"
$bBot[1,2] = "Free"
$bBot[1,3] = "Free"
$bBot[1,4] = "Free"
$bBot[2,0] = "Free"
$bBot[2,1] = "CS309"
$bBot[2,2] = "CS309"
$bBot[2,3] = "CS304"
$bBot[2,4] = "Free"
$bBot[3,0] = "Free"
$bBot[3,1] = "Free"
$bBot[3,2] = "Free"
$bBot[3,3] = "Free"
$bBot[3,4] = "CS304 (Tutorial)"
$wsB.Range("B5:F8").Value = $bBot

# ---- Course_Summary table rewrite ----
$wsC = $wb.Worksheets.Item("Course_Summary")
$cs = New-Object 'object[,]' 33,10
$cs[0,0] = "Course Code"
$cs[0,1] = "Course Name"
$cs[0,2] = "Course Type"
$cs[0,3] = "Branch Specificity"
$cs[0,4] = "LTPSC"
$cs[0,5] = "Lectures/Week"
$cs[0,6] = "Tutorials/Week"
$cs[0,7] = "Total Credits"
$cs[0,8] = "Instructor"
$cs[0,9] = "Department"
$cs[1,0] = "CS309"
$cs[1,1] = "Statistics for CS"
$cs[1,2] = "Core"
$cs[1,3] = "Department: CSE"
$cs[1,4] = "3-1-0-0-4"
$cs[1,5] = 3
$cs[1,6] = 1
$cs[1,7] = 4
$cs[1,8] = "Sunil C K, Pavan"
$cs[1,9] = "CSE"
$cs[2,0] = "CS303"
$cs[2,1] = "Computer Networks"
$cs[2,2] = "Core"
$cs[2,3] = "Department: CSE"
$cs[2,4] = "3-1-2-0-5"
$cs[2,5] = 3
$cs[2,6] = 1
$cs[2,7] = 5
$cs[2,8] = "Animesh Roy, Dibyajyothi"
$cs[2,9] = "CSE"
$cs[3,0] = "CS304"
$cs[3,1] = "Artificial Intelligence"
$cs[3,2] = "Core"
$cs[3,3] = "Department: CSE"
$cs[3,4] = "3-1-0-0-4"
$cs[3,5] = 3
$cs[3,6] = 1
$cs[3,7] = 4
$cs[3,8] = "Krishendu, Girish"
$cs[3,9] = "CSE"
$cs[4,0] = "HS101"
$cs[4,1] = "Environmental Studies"
$cs[4,2] = "Core"
$cs[4,3] = "Department: CSE"
$cs[4,4] = "0-0-0-8-2"
$cs[4,5] = 0
$cs[4,6] = 0
$cs[4,7] = 2
$cs[4,8] = "TBD"
$cs[4,9] = "CSE"
$cs[5,0] = "CS463"
$cs[5,1] = "Parallel computing"
$cs[5,2] = "Elective"
$cs[5,3] = "Common for All Branches"
$cs[5,4] = "3-1-0-0-4"
$cs[5,5] = 3
$cs[5,6] = 1
$cs[5,7] = 4
$cs[5,8] = "Pramod"
$cs[5,9] = "CSE"
$cs[6,0] = "CS308"
$cs[6,1] = "Compiler Design"
$cs[6,2] = "Elective"
$cs[6,3] = "Common for All Branches"
$cs[6,4] = "3-1-0-0-4"
$cs[6,5] = 3
$cs[6,6] = 1
$cs[6,7] = 4
$cs[6,8] = "Pavan Kumar"
$cs[6,9] = "CSE"
$cs[7,0] = "DS301"
$cs[7,1] = "Graphs and Social Networks"
$cs[7,2] = "Elective"
$cs[7,3] = "Common for All Branches"
$cs[7,4] = "3-1-0-0-4"
$cs[7,5] = 3
$cs[7,6] = 1
$cs[7,7] = 4
$cs[7,8] = "Utkarsh K"
$cs[7,9] = "CSE"
$cs[8,0] = "CS366"
$cs[8,1] = "Advanced Algorithms"
$cs[8,2] = "Elective"
$cs[8,3] = "Common for All Branches"
$cs[8,4] = "3-1-0-0-4"
$cs[8,5] = 3
$cs[8,6] = 1
$cs[8,7] = 4
$cs[8,8] = "Suvadip"
$cs[8,9] = "CSE"
$cs[9,0] = "DS359"
$cs[9,1] = "Full Stack Development"
$cs[9,2] = "Elective"
$cs[9,3] = "Common for All Branches"
$cs[9,4] = "3-1-0-0-4"
$cs[9,5] = 3
$cs[9,6] = 1
$cs[9,7] = 4
$cs[9,8] = "Manjunath K V"
$cs[9,9] = "CSE"
$cs[10,0] = "EC355"
$cs[10,1] = "Internet of Things"
$cs[10,2] = "Elective"
$cs[10,3] = "Common for All Branches"
$cs[10,4] = "3-1-0-0-4"
$cs[10,5] = 3
$cs[10,6] = 1
$cs[10,7] = 4
$cs[10,8] = "Prakash Pawar"
$cs[10,9] = "CSE"
$cs[11,0] = "EC364"
$cs[11,1] = "Semiconductor Devices"
$cs[11,2] = "Elective"
$cs[11,3] = "Common for All Branches"
$cs[11,4] = "3-1-0-0-4"
$cs[11,5] = 3
$cs[11,6] = 1
$cs[11,7] = 4
$cs[11,8] = "Pankaj Kumar"
$cs[11,9] = "CSE"
$cs[12,0] = "CS352"
$cs[12,1] = "Cryptography & Security"
$cs[12,2] = "Elective"
$cs[12,3] = "Common for All Branches"
$cs[12,4] = "3-1-0-0-4"
$cs[12,5] = 3
$cs[12,6] = 1
$cs[12,7] = 4
$cs[12,8] = "Rajendra Hegadi"
$cs[12,9] = "CSE"
$cs[13,0] = "ASD352"
$cs[13,1] = "User Interaction"
$cs[13,2] = "Elective"
$cs[13,3] = "Common for All Branches"
$cs[13,4] = "3-1-0-0-4"
$cs[13,5] = 3
$cs[13,6] = 1
$cs[13,7] = 4
$cs[13,8] = "Sandesh P"
$cs[13,9] = "CSE"
$cs[14,0] = "EC365"
$cs[14,1] = "AI in Biomedical Engineering"
$cs[14,2] = "Elective"
$cs[14,3] = "Common for All Branches"
$cs[14,4] = "3-1-0-0-4"
$cs[14,5] = 3
$cs[14,6] = 1
$cs[14,7] = 4
$cs[14,8] = "Sibasankar Padhy"
$cs[14,9] = "CSE"
$cs[15,0] = "CS463"
$cs[15,1] = "Parallel computing"
$cs[15,2] = "Elective"
$cs[15,3] = "Common for All Branches"
$cs[15,4] = "3-1-0-0-4"
$cs[15,5] = 3
$cs[15,6] = 1
$cs[15,7] = 4
$cs[15,8] = "Pramod"
$cs[15,9] = "DSAI"
$cs[16,0] = "CS308"
$cs[16,1] = "Compiler Design"
$cs[16,2] = "Elective"
$cs[16,3] = "Common for All Branches"
$cs[16,4] = "3-1-0-0-4"
$cs[16,5] = 3
$cs[16,6] = 1
$cs[16,7] = 4
$cs[16,8] = "Pavan Kumar"
$cs[16,9] = "DSAI"
$cs[17,0] = "DS301"
$cs[17,1] = "Graphs and Social Networks"
$cs[17,2] = "Elective"
$cs[17,3] = "Common for All Branches"
$cs[17,4] = "3-1-0-0-4"
$cs[17,5] = 3
$cs[17,6] = 1
$cs[17,7] = 4
$cs[17,8] = "Utkarsh K"
$cs[17,9] = "DSAI"
$cs[18,0] = "CS366"
$cs[18,1] = "Advanced Algorithms"
$cs[18,2] = "Elective"
$cs[18,3] = "Common for All Branches"
$cs[18,4] = "3-1-0-0-4"
$cs[18,5] = 3
$cs[18,6] = 1
$cs[18,7] = 4
$cs[18,8] = "Suvadip"
$cs[18,9] = "DSAI"
$cs[19,0] = "DS359"
$cs[19,1] = "Full Stack Development"
$cs[19,2] = "Elective"
$cs[19,3] = "Common for All Branches"
$cs[19,4] = "3-1-0-0-4"
$cs[19,5] = 3
$cs[19,6] = 1
$cs[19,7] = 4
$cs[19,8] = "Manjunath K V"
$cs[19,9] = "DSAI"
$cs[20,0] = "EC355"
$cs[20,1] = "Internet of Things"
$cs[20,2] = "Elective"
$cs[20,3] = "Common for All Branches"
$cs[20,4] = "3-1-0-0-4"
$cs[20,5] = 3
$cs[20,6] = 1
$cs[20,7] = 4
$cs[20,8] = "Prakash Pawar"
$cs[20,9] = "DSAI"
$cs[21,0] = "EC364"
$cs[21,1] = "Semiconductor Devices"
$cs[21,2] = "Elective"
$cs[21,3] = "Common for All Branches"
$cs[21,4] = "3-1-0-0-4"
$cs[21,5] = 3
$cs[21,6] = 1
$cs[21,7] = 4
$cs[21,8] = "Pankaj Kumar"
$cs[21,9] = "DSAI"
$cs[22,0] = "CS352"
$cs[22,1] = "Cryptography & Security"
$cs[22,2] = "Elective"
$cs[22,3] = "Common for All Branches"
$cs[22,4] = "3-1-0-0-4"
$cs[22,5] = 3
$cs[22,6] = 1
$cs[22,7] = 4
$cs[22,8] = "Rajendra Hegadi"
$cs[22,9] = "DSAI"
$cs[23,0] = "ASD352"
$cs[23,1] = "User Interaction"
$cs[23,2] = "Elective"
$cs[23,3] = "Common for All Branches"
$cs[23,4] = "3-1-0-0-4"
$cs[23,5] = 3
$cs[23,6] = 1
$cs[23,7] = 4
$cs[23,8] = "Sandesh P"
$cs[23,9] = "DSAI"
$cs[24,0] = "EC365"
$cs[24,1] = "AI in Biomedical Engineering"
$cs[24,2] = "Elective"
$cs[24,3] = "Common for All Branches"
$cs[24,4] = "3-1-0-0-4"
$cs[24,5] = 3
$cs[24,6] = 1
$cs[24,7] = 4
$cs[24,8] = "Sibasankar Padhy"
$cs[24,9] = "DSAI"
$cs[25,0] = "HS301"
$cs[25,1] = "HSS/IE"
$cs[25,2] = "Elective"
$cs[25,3] = "Common for All Branches"
$cs[25,4] = "3-0-0-0-3"
$cs[25,5] = 3
$cs[25,6] = 0
$cs[25,7] = 3
$cs[25,8] = "TBD"
$cs[25,9] = "ECE"
$cs[26,0] = "EC251"
$cs[26,1] = "Electronics Systems"
$cs[26,2] = "Elective"
$cs[26,3] = "Common for All Branches"
$cs[26,4] = "2-0-0-0-2"
$cs[26,5] = 2
$cs[26,6] = 0
$cs[26,7] = 2
$cs[26,8] = "Pankaj"
$cs[26,9] = "ECE"
$cs[27,0] = "EC252"
$cs[27,1] = "Introduction of Communication"
$cs[27,2] = "Elective"
$cs[27,3] = "Common for All Branches"
$cs[27,4] = "2-0-0-0-2"
$cs[27,5] = 2
$cs[27,6] = 0
$cs[27,7] = 2
$cs[27,8] = "Rajesh Kumar"
$cs[27,9] = "ECE"
$cs[28,0] = "EC253"
$cs[28,1] = "Introduction to AI"
$cs[28,2] = "Elective"
$cs[28,3] = "Common for All Branches"
$cs[28,4] = "2-0-2-4-2"
$cs[28,5] = 2
$cs[28,6] = 0
$cs[28,7] = 2
$cs[28,8] = "Sibasankar Padhy"
$cs[28,9] = "ECE"
$cs[29,0] = "EC254"
$cs[29,1] = "Electronic Systems"
$cs[29,2] = "Elective"
$cs[29,3] = "Common for All Branches"
$cs[29,4] = "2-0-0-0-2"
$cs[29,5] = 2
$cs[29,6] = 0
$cs[29,7] = 2
$cs[29,8] = "Mallikarjun Kande"
$cs[29,9] = "ECE"
$cs[30,0] = "CS162"
$cs[30,1] = "Data Science with Python"
$cs[30,2] = "Elective"
$cs[30,3] = "Common for All Branches"
$cs[30,4] = "2-0-0-0-2"
$cs[30,5] = 2
$cs[30,6] = 0
$cs[30,7] = 2
$cs[30,8] = "Abdul Wahid"
$cs[30,9] = "ECE"
$cs[31,0] = "DE352"
$cs[31,1] = "User Interaction (Minor)"
$cs[31,2] = "Elective"
$cs[31,3] = "Common for All Branches"
$cs[31,4] = "3-1-0-0-4"
$cs[31,5] = 3
$cs[31,6] = 1
$cs[31,7] = 4
$cs[31,8] = "Sandesh P"
$cs[31,9] = "ECE"
$cs[32,0] = "CS251"
$cs[32,1] = "2D Computer Graphics"
$cs[32,2] = "Elective"
$cs[32,3] = "Common for All Branches"
$cs[32,4] = "2-0-0-0-2"
$cs[32,5] = 2
$cs[32,6] = 0
$cs[32,7] = 2
$cs[32,8] = "Vivekraj"
$cs[32,9] = "ECE"
$wsC.Range("A1:J33").Value = $cs

# Remove leftover rows 34-41 from the old 41-row table
$wsC.Range("A34:J41").EntireRow.Delete()

# Apply header style (bold/centered/bordered, matching style index 1) to new Department header cell
$wsC.Range("J1").Style = $wsC.Range("I1").Style

$wsC.Range("A1:J33").Columns.AutoFit() | Out-Null
